$d = $word.ActiveDocument

# Locate the bullet paragraph "Use responsive web design frameworks like w3.css or bootstrap"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Use responsive web design frameworks like w3.css or bootstrap*") {
        $target = $p
        break
    }
}
if (-not $target) {
    throw "Could not locate the 'Use responsive web design frameworks...' paragraph"
}

# The paragraph that currently follows it (the empty bookmark paragraph) - inserting a new
# paragraph break just before it gives us a brand-new, plain paragraph (no list/pStyle
# inheritance) sitting between the two.
$followingPara = $target.Next()
[void]$followingPara.Range.InsertParagraphBefore()

# Re-resolve the freshly created (still empty) paragraph.
$newPara = $target.Next()
$newRange = $newPara.Range

# Build the paragraph content ("HTML " + "Computercode" split across two runs with a
# spell-check proofErr pair around "Computercode", all underlined) as raw OOXML and inject
# it via InsertXML, which replaces the (empty) contents of $newRange in place.
$fragXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">HTML </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>Computercode</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$newRange.InsertXML($fragXml)
